$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "done" (column M) checklist cells with an "X" for the
# user stories that are now finished.
$doneRows = @(4, 5, 7, 8, 9, 11, 12, 13, 14, 15, 19)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 13).Value = "X"
}

# Mark the "being developed" (column J) checklist cells with an "X"
# for the user stories currently in progress.
$inProgressRows = @(17, 20, 22, 23)
foreach ($r in $inProgressRows) {
    $ws.Cells.Item($r, 10).Value = "X"
}

# Refresh the sheet view: scroll/zoom and selection as left by the author.
$ws.Activate() | Out-Null
$ws.Range("M20").Select() | Out-Null
$excel.ActiveWindow.Zoom = 60
